$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new daily rows (date serial, nuovi pos. = 0, somma mobile 7gg. = 0, per 100mila ab. = 0)
$newRows = @(
    @{Row=252; Date=44326},
    @{Row=253; Date=44327},
    @{Row=254; Date=44328},
    @{Row=255; Date=44329}
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Copy formatting (style) from the cell directly above, same as extending the date column
    $ws.Cells.Item($row - 1, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
}

$excel.CutCopyMode = $false
